$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; existing rows 9..103 shift down to 10..104.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44552
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 100112052
$ws.Cells.Item(9, 7).Value = "Albahaca"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 125
$ws.Cells.Item(9, 11).Value = 4000
$ws.Cells.Item(9, 12).Value = 4500
$ws.Cells.Item(9, 13).Value = 4260
$ws.Cells.Item(9, 14).Value = "$/docena de matas"
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 710
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = "Hortaliza"
